$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 5261.222
$ws.Cells.Item(88, 10).Value = 5641.8335
$ws.Cells.Item(88, 12).Value = 5641.8335
$ws.Cells.Item(88, 14).Value = -6453.8335
$ws.Cells.Item(91, 8).Value = 5261.222
$ws.Cells.Item(91, 10).Value = 5641.8335
$ws.Cells.Item(91, 12).Value = 5641.8335
$ws.Cells.Item(91, 14).Value = -8449.833500000001
$ws.Cells.Item(129, 8).Value = 107256.58
$ws.Cells.Item(129, 9).Value = 168031.33
$ws.Cells.Item(129, 11).Value = 504093.99
$ws.Cells.Item(129, 13).Value = -499093.99
$ws.Cells.Item(132, 8).Value = 16833.422
$ws.Cells.Item(132, 9).Value = 2764.0645
$ws.Cells.Item(132, 11).Value = 8292.193499999999
$ws.Cells.Item(132, 13).Value = -5762.193499999999
$ws.Cells.Item(135, 8).Value = 12501413
$ws.Cells.Item(135, 9).Value = 14287071
$ws.Cells.Item(135, 10).Value = 1804.4
$ws.Cells.Item(135, 11).Value = 128583639
$ws.Cells.Item(135, 12).Value = 16239.6
$ws.Cells.Item(135, 13).Value = -128581104
$ws.Cells.Item(135, 14).Value = -21309.6
$ws.Cells.Item(141, 8).Value = 4852.516
$ws.Cells.Item(141, 9).Value = 2158.1428
$ws.Cells.Item(141, 11).Value = 6474.428400000001
$ws.Cells.Item(141, 13).Value = -1294.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 2801
$ws.Cells.Item(14, 10).Value = 3568
$ws.Cells.Item(14, 12).Value = 3568
$ws.Cells.Item(14, 14).Value = -3918
$ws.Cells.Item(32, 8).Value = 21744956
$ws.Cells.Item(32, 9).Value = 21744956
$ws.Cells.Item(32, 11).Value = 21744956
$ws.Cells.Item(32, 13).Value = -21744669
$ws.Cells.Item(63, 8).Value = 9379
$ws.Cells.Item(63, 9).Value = 2298.3333
$ws.Cells.Item(63, 11).Value = 2298.3333
$ws.Cells.Item(63, 13).Value = -1612.3333
$ws.Cells.Item(66, 8).Value = 9379
$ws.Cells.Item(66, 9).Value = 2298.3333
$ws.Cells.Item(66, 11).Value = 11491.6665
$ws.Cells.Item(66, 13).Value = -8059.666499999999
$ws.Cells.Item(74, 8).Value = 1788.8928
$ws.Cells.Item(74, 9).Value = 1788.8928
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 1788.8928
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -914.8928000000001
$ws.Cells.Item(74, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 1788.8928
$ws.Cells.Item(77, 9).Value = 1788.8928
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 8944.464
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -4576.464
$ws.Cells.Item(77, 14).ClearContents()
$ws.Cells.Item(92, 8).Value = 51912.5
$ws.Cells.Item(92, 10).Value = 51912.5
$ws.Cells.Item(92, 12).Value = 51912.5
$ws.Cells.Item(92, 14).Value = -56904.5
$ws.Cells.Item(132, 8).Value = 1552.5682
$ws.Cells.Item(132, 9).Value = 1566.9762
$ws.Cells.Item(132, 11).Value = 4700.9286
$ws.Cells.Item(132, 13).Value = -2170.9286

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 13122.5
$ws.Cells.Item(107, 9).Value = 9163.333000000001
$ws.Cells.Item(107, 10).Value = 25000
$ws.Cells.Item(107, 11).Value = 9163.333000000001
$ws.Cells.Item(107, 12).Value = 25000
$ws.Cells.Item(107, 13).Value = -7243.333000000001
$ws.Cells.Item(107, 14).Value = -28840
$ws.Cells.Item(134, 8).Value = 2191.1785
$ws.Cells.Item(134, 9).Value = 1821.3077
$ws.Cells.Item(134, 10).Value = 6999.5
$ws.Cells.Item(134, 11).Value = 5463.9231
$ws.Cells.Item(134, 12).Value = 20998.5
$ws.Cells.Item(134, 13).Value = -2928.9231
$ws.Cells.Item(134, 14).Value = -26068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 149.21428
$ws.Cells.Item(7, 10).Value = 224.25
$ws.Cells.Item(7, 12).Value = 224.25
$ws.Cells.Item(7, 14).Value = -450.25
$ws.Cells.Item(31, 8).Value = 1815.5264
$ws.Cells.Item(31, 9).Value = 1766.8276
$ws.Cells.Item(31, 10).Value = 1972.4445
$ws.Cells.Item(31, 11).Value = 1766.8276
$ws.Cells.Item(31, 12).Value = 1972.4445
$ws.Cells.Item(31, 13).Value = -1471.8276
$ws.Cells.Item(31, 14).Value = -2562.4445
$ws.Cells.Item(34, 8).Value = 1815.5264
$ws.Cells.Item(34, 9).Value = 1766.8276
$ws.Cells.Item(34, 10).Value = 1972.4445
$ws.Cells.Item(34, 11).Value = 1766.8276
$ws.Cells.Item(34, 12).Value = 1972.4445
$ws.Cells.Item(34, 13).Value = -1564.8276
$ws.Cells.Item(34, 14).Value = -2376.4445
$ws.Cells.Item(58, 8).Value = 1149.125
$ws.Cells.Item(58, 9).Value = 768
$ws.Cells.Item(58, 11).Value = 768
$ws.Cells.Item(58, 13).Value = -565
$ws.Cells.Item(136, 8).Value = 1149.125
$ws.Cells.Item(136, 9).Value = 768
$ws.Cells.Item(136, 11).Value = 2304
$ws.Cells.Item(136, 13).Value = 246

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 32419636
$ws.Cells.Item(4, 9).Value = 35893064
$ws.Cells.Item(4, 11).Value = 107679192
$ws.Cells.Item(4, 13).Value = -107679080
$ws.Cells.Item(7, 8).Value = 1122.25
$ws.Cells.Item(7, 9).Value = 171.33333
$ws.Cells.Item(7, 11).Value = 513.99999
$ws.Cells.Item(7, 13).Value = -401.99999
$ws.Cells.Item(11, 8).Value = 3285.889
$ws.Cells.Item(11, 9).Value = 4818.25
$ws.Cells.Item(11, 10).Value = 221.16667
$ws.Cells.Item(11, 11).Value = 14454.75
$ws.Cells.Item(11, 12).Value = 663.50001
$ws.Cells.Item(11, 13).Value = -14314.75
$ws.Cells.Item(11, 14).Value = -943.50001
$ws.Cells.Item(59, 8).Value = 214816.33
$ws.Cells.Item(59, 10).Value = 319449
$ws.Cells.Item(59, 12).Value = 958347
$ws.Cells.Item(59, 14).Value = -959427
$ws.Cells.Item(107, 8).Value = 2950
$ws.Cells.Item(107, 10).Value = 4250
$ws.Cells.Item(107, 12).Value = 12750
$ws.Cells.Item(107, 14).Value = -16590
$ws.Cells.Item(125, 8).Value = 6035
$ws.Cells.Item(125, 9).Value = 5552.5
$ws.Cells.Item(125, 11).Value = 16657.5
$ws.Cells.Item(125, 13).Value = -11737.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 8191.5713
$ws.Cells.Item(5, 9).Value = 8191.5713
$ws.Cells.Item(5, 11).Value = 8191.5713
$ws.Cells.Item(5, 13).Value = -8079.5713
$ws.Cells.Item(14, 8).Value = 746392.8
$ws.Cells.Item(14, 9).Value = 1215642.9
$ws.Cells.Item(14, 10).Value = 277142.72
$ws.Cells.Item(14, 11).Value = 1215642.9
$ws.Cells.Item(14, 12).Value = 277142.72
$ws.Cells.Item(14, 13).Value = -1215474.9
$ws.Cells.Item(14, 14).Value = -277478.72
$ws.Cells.Item(19, 8).Value = 698.6667
$ws.Cells.Item(19, 9).Value = 550
$ws.Cells.Item(19, 10).Value = 996
$ws.Cells.Item(19, 11).Value = 550
$ws.Cells.Item(19, 12).Value = 996
$ws.Cells.Item(19, 13).Value = -262
$ws.Cells.Item(19, 14).Value = -1572
$ws.Cells.Item(22, 8).Value = 450
$ws.Cells.Item(22, 9).Value = 450
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 450
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 79
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 6999.393
$ws.Cells.Item(80, 9).Value = 5745.353
$ws.Cells.Item(80, 11).Value = 5745.353
$ws.Cells.Item(80, 13).Value = -4747.353
$ws.Cells.Item(83, 8).Value = 6999.393
$ws.Cells.Item(83, 9).Value = 5745.353
$ws.Cells.Item(83, 11).Value = 28726.765
$ws.Cells.Item(83, 13).Value = -23734.765
$ws.Cells.Item(94, 8).Value = 30666.6
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 30666.6
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 30666.6
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -32018.6
$ws.Cells.Item(97, 8).Value = 416.55554
$ws.Cells.Item(97, 9).Value = 379.2143
$ws.Cells.Item(97, 10).Value = 547.25
$ws.Cells.Item(97, 11).Value = 379.2143
$ws.Cells.Item(97, 12).Value = 547.25
$ws.Cells.Item(97, 13).Value = 116.7857
$ws.Cells.Item(97, 14).Value = -1539.25
$ws.Cells.Item(132, 8).Value = 2402.4
$ws.Cells.Item(132, 9).Value = 2503
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 7509
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -4979
$ws.Cells.Item(132, 14).Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 12499.25
$ws.Cells.Item(93, 9).Value = 9999
$ws.Cells.Item(93, 10).Value = 20000
$ws.Cells.Item(93, 11).Value = 9999
$ws.Cells.Item(93, 12).Value = 20000
$ws.Cells.Item(93, 13).Value = -8751
$ws.Cells.Item(93, 14).Value = -22496
$ws.Cells.Item(100, 8).Value = 6099.4736
$ws.Cells.Item(100, 9).Value = 3987.3333
$ws.Cells.Item(100, 11).Value = 3987.3333
$ws.Cells.Item(100, 13).Value = -3446.3333
$ws.Cells.Item(122, 8).Value = 5494.5557
$ws.Cells.Item(122, 9).Value = 2744.8
$ws.Cells.Item(122, 11).Value = 8234.400000000001
$ws.Cells.Item(122, 13).Value = -5784.400000000001
$ws.Cells.Item(132, 8).Value = 4576.3125
$ws.Cells.Item(132, 9).Value = 2070.6155
$ws.Cells.Item(132, 10).Value = 15434.333
$ws.Cells.Item(132, 11).Value = 6211.8465
$ws.Cells.Item(132, 12).Value = 46302.999
$ws.Cells.Item(132, 13).Value = -3681.8465
$ws.Cells.Item(132, 14).Value = -51362.999
$ws.Cells.Item(136, 8).Value = 3246.0386
$ws.Cells.Item(136, 9).Value = 3321.6086
$ws.Cells.Item(136, 11).Value = 9964.825800000001
$ws.Cells.Item(136, 13).Value = -7414.825800000001
$ws.Cells.Item(139, 8).Value = 944266.5
$ws.Cells.Item(139, 10).Value = 1085606.5
$ws.Cells.Item(139, 12).Value = 1085606.5
$ws.Cells.Item(139, 14).Value = -1095886.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 46000
$ws.Cells.Item(16, 10).Value = 46000
$ws.Cells.Item(16, 12).Value = 46000
$ws.Cells.Item(16, 14).Value = -46584
$ws.Cells.Item(121, 8).Value = 14997.5
$ws.Cells.Item(121, 10).Value = 14997.5
$ws.Cells.Item(121, 12).Value = 14997.5
$ws.Cells.Item(121, 14).Value = -18491.5
$ws.Cells.Item(122, 8).Value = 2320.5
$ws.Cells.Item(122, 9).Value = 2320.5
$ws.Cells.Item(122, 11).Value = 6961.5
$ws.Cells.Item(122, 13).Value = -4511.5
$ws.Cells.Item(126, 8).Value = 1536.08
$ws.Cells.Item(126, 9).Value = 1405.5238
$ws.Cells.Item(126, 11).Value = 4216.5714
$ws.Cells.Item(126, 13).Value = -1746.5714
